$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before row 302, shifting existing rows 302-321 down to 304-323.
$ws.Rows("302:303").Insert()

# Row 302 - new record
$ws.Cells.Item(302, 1).Value = 5
$ws.Cells.Item(302, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(302, 3).Value = "Maule"
$ws.Cells.Item(302, 4).Value = 44516
$ws.Cells.Item(302, 5).Value = 7
$ws.Cells.Item(302, 6).Value = 100114001
$ws.Cells.Item(302, 7).Value = "Papa"
$ws.Cells.Item(302, 8).Value = "Asterix"
$ws.Cells.Item(302, 9).Value = "1a nueva(o)"
$ws.Cells.Item(302, 10).Value = 1600
$ws.Cells.Item(302, 11).Value = 10000
$ws.Cells.Item(302, 12).Value = 10000
$ws.Cells.Item(302, 13).Value = 10000
$ws.Cells.Item(302, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(302, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(302, 16).Value = 400
$ws.Cells.Item(302, 17).Value = 25
$ws.Cells.Item(302, 18).Value = "Hortaliza"

# Row 303 - new record
$ws.Cells.Item(303, 1).Value = 5
$ws.Cells.Item(303, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(303, 3).Value = "Maule"
$ws.Cells.Item(303, 4).Value = 44516
$ws.Cells.Item(303, 5).Value = 7
$ws.Cells.Item(303, 6).Value = 100114001
$ws.Cells.Item(303, 7).Value = "Papa"
$ws.Cells.Item(303, 8).Value = "Rodeo"
$ws.Cells.Item(303, 9).Value = "1a nueva(o)"
$ws.Cells.Item(303, 10).Value = 1300
$ws.Cells.Item(303, 11).Value = 10000
$ws.Cells.Item(303, 12).Value = 10000
$ws.Cells.Item(303, 13).Value = 10000
$ws.Cells.Item(303, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(303, 15).Value = "Región del Maule"
$ws.Cells.Item(303, 16).Value = 400
$ws.Cells.Item(303, 17).Value = 25
$ws.Cells.Item(303, 18).Value = "Hortaliza"
